# Commit: Tue, Mar 24, 2020  2:05:27 AM
#
# The deck's "Integral" design theme (ppt/theme/theme1.xml, the theme used
# by the one SlideMaster / all slides) is swapped for the stock "Office
# Theme" color palette that already ships - unused by any slide - as
# ppt/theme/theme2.xml (the Notes Master's theme). The table on slide 16
# is re-pointed from the deck's custom "Table_0" table style to the
# built-in "Medium Style 2 - Accent 1" style.

$p = $ppt.ActivePresentation

# --- 1. Slide 16's table: switch to the built-in table style -------------
$slide  = $p.Slides.Item(16)
$shape  = $slide.Shapes.Item(3)          # the single p:graphicFrame/table
$table  = $shape.Table
$table.ApplyStyle("{76418BAA-083E-4F7B-AA71-47B4B230490F}")

# --- 2. Re-theme the deck: Integral -> Office Theme colors ----------------
# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as the BGR-packed ints PowerPoint's RGB()/.RGB wants.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
